# ---------------------------------------------------------------------------
# "Actualizando PIB,VAB y resultado fiscal"
#
# 1. Sheet "Producto" (PIB): refresh the Q-on-Q PIB DESEST series in column C
#    (rows 2-66), refresh D66, and append a new data point (row 67: 2020-II).
# 2. Sheet "VAB": refresh the "Valor" / "Variacion anual" columns (B:C,
#    rows 2-18) - the VAB + resultado fiscal figures.
# 3. Minor formatting touch-up (center alignment) on the Pobreza sheets'
#    header rows, matching the workbook author's pass.
# 4. View-state: the author had scrolled to/selected a different part of
#    "Producto" and left the "VAB" tab active/selected when saving.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet "Producto": refresh column C (rows 2-66) + D66, add row 67 ---
$wsProducto = $wb.Worksheets.Item("Producto")

$wsProducto.Range("C2").Value = 475540.07912460313
$wsProducto.Range("C3").Value = 470121.1274420924
$wsProducto.Range("C4").Value = 493832.66195269296
$wsProducto.Range("C5").Value = 500966.91203686345
$wsProducto.Range("C6").Value = 515472.67133145465
$wsProducto.Range("C7").Value = 526263.5798129034
$wsProducto.Range("C8").Value = 530051.5527618445
$wsProducto.Range("C9").Value = 540435.966107454
$wsProducto.Range("C10").Value = 554527.005412419
$wsProducto.Range("C11").Value = 561248.738649315
$wsProducto.Range("C12").Value = 576963.188307025
$wsProducto.Range("C13").Value = 589458.6845141673
$wsProducto.Range("C14").Value = 603178.7291797571
$wsProducto.Range("C15").Value = 616364.1681524651
$wsProducto.Range("C16").Value = 624543.4795623703
$wsProducto.Range("C17").Value = 643683.6336897694
$wsProducto.Range("C18").Value = 649543.926967607
$wsProducto.Range("C19").Value = 653526.3466845228
$wsProducto.Range("C20").Value = 658399.4477377966
$wsProducto.Range("C21").Value = 627234.9175749234
$wsProducto.Range("C22").Value = 604302.8264205777
$wsProducto.Range("C23").Value = 591301.3909057183
$wsProducto.Range("C24").Value = 614138.8940492504
$wsProducto.Range("C25").Value = 625748.3942759503
$wsProducto.Range("C26").Value = 644736.4091832255
$wsProducto.Range("C27").Value = 673861.030548019
$wsProducto.Range("C28").Value = 677043.3951807589
$wsProducto.Range("C29").Value = 686453.8828551905
$wsProducto.Range("C30").Value = 702907.2608763932
$wsProducto.Range("C31").Value = 709431.716286141
$wsProducto.Range("C32").Value = 715005.9318030047
$wsProducto.Range("C33").Value = 715781.4799168685
$wsProducto.Range("C34").Value = 708078.0613420046
$wsProducto.Range("C35").Value = 683459.320640823
$wsProducto.Range("C36").Value = 705124.7536755278
$wsProducto.Range("C37").Value = 717281.8221774418
$wsProducto.Range("C38").Value = 717239.1536081125
$wsProducto.Range("C39").Value = 720685.010983557
$wsProducto.Range("C40").Value = 725384.147289164
$wsProducto.Range("C41").Value = 718320.1093304262
$wsProducto.Range("C42").Value = 707643.7549808683
$wsProducto.Range("C43").Value = 703081.7383881108
$wsProducto.Range("C44").Value = 697440.750144984
$wsProducto.Range("C45").Value = 701057.940339497
$wsProducto.Range("C46").Value = 711623.6110289416
$wsProducto.Range("C47").Value = 727716.4679055146
$wsProducto.Range("C48").Value = 727288.6519749603
$wsProducto.Range("C49").Value = 719319.8556427438
$wsProducto.Range("C50").Value = 713343.5535243312
$wsProducto.Range("C51").Value = 700876.9889477105
$wsProducto.Range("C52").Value = 703503.3059038705
$wsProducto.Range("C53").Value = 708187.5460147362
$wsProducto.Range("C54").Value = 715541.1243024282
$wsProducto.Range("C55").Value = 720864.5454440481
$wsProducto.Range("C56").Value = 730409.7080909061
$wsProducto.Range("C57").Value = 738744.4132139125
$wsProducto.Range("C58").Value = 738063.5060731241
$wsProducto.Range("C59").Value = 700381.8126628962
$wsProducto.Range("C60").Value = 699912.3783340706
$wsProducto.Range("C61").Value = 692664.2613444821
$wsProducto.Range("C62").Value = 693888.025796307
$wsProducto.Range("C63").Value = 690911.4545944538
$wsProducto.Range("C64").Value = 696760.174579235
$wsProducto.Range("C65").Value = 690350.1478867279
$wsProducto.Range("C66").Value = 661654.9489224798

$wsProducto.Range("D66").Value = 630971.4937514844

# New row 67: 2020 - II quarter
$wsProducto.Range("A67").Value = 2020
$wsProducto.Range("B67").Value = "II"
$wsProducto.Range("C67").Value = 554316.13938229729
$wsProducto.Range("D67").Value = 606992.80173888081
$wsProducto.Range("A67:C67").HorizontalAlignment = -4108
$wsProducto.Range("D67").Style = "Normal"


# --- Sheet "VAB": refresh columns B (Valor) and C (Variacion anual), rows 2-18 ---
$wsVAB = $wb.Worksheets.Item("VAB")

$wsVAB.Range("B2").Value = 512359.5766835184
$wsVAB.Range("C2").Value = -19.36708585118373
$wsVAB.Range("B3").Value = 92004.64506888161
$wsVAB.Range("C3").Value = -10.662776074013092
$wsVAB.Range("B4").Value = 1762.0229144674997
$wsVAB.Range("C4").Value = -14.01278029575922
$wsVAB.Range("B5").Value = 17786.684208885996
$wsVAB.Range("C5").Value = -18.257336008282532
$wsVAB.Range("B6").Value = 87855.64006166624
$wsVAB.Range("C6").Value = -20.78211950432992
$wsVAB.Range("B7").Value = 11444.28144229284
$wsVAB.Range("C7").Value = -3.322617578802578
$wsVAB.Range("B8").Value = 10168.960336140055
$wsVAB.Range("C8").Value = -52.06171237138501
$wsVAB.Range("B9").Value = 71161.53514765836
$wsVAB.Range("C9").Value = -16.915809495210297
$wsVAB.Range("B10").Value = 2846.6104571211044
$wsVAB.Range("C10").Value = -73.35874308992452
$wsVAB.Range("B11").Value = 45942.355822993144
$wsVAB.Range("C11").Value = -22.455757918651422
$wsVAB.Range("B12").Value = 24870.172648461972
$wsVAB.Range("C12").Value = -1.2002304258781127
$wsVAB.Range("B13").Value = 64522.940872450235
$wsVAB.Range("C13").Value = -14.315387598774443
$wsVAB.Range("B14").Value = 28996.86758796262
$wsVAB.Range("C14").Value = -12.799845828932533
$wsVAB.Range("B15").Value = 25007.616698329046
$wsVAB.Range("C15").Value = -9.209401485796686
$wsVAB.Range("B16").Value = 19525.68352757101
$wsVAB.Range("C16").Value = -23.468530609699613
$wsVAB.Range("B17").Value = 5830.519916386826
$wsVAB.Range("C17").Value = -67.71281457157797
$wsVAB.Range("B18").Value = 2633.0399722498255
$wsVAB.Range("C18").Value = -38.02133999999998

# --- Formatting touch-up: center the header row (B1:E1) on the Pobreza sheets ---
$wsPobrezaAglo = $wb.Worksheets.Item("Pobreza-Aglo")
$wsPobrezaAglo.Range("B1:E1").HorizontalAlignment = -4108

$wsPobrezaRegiones = $wb.Worksheets.Item("Pobreza regiones")
$wsPobrezaRegiones.Range("B1:E1").HorizontalAlignment = -4108

# --- View state ---
# "Producto": leave the view scrolled to/selected around F63 (bottom of data)
$wsProducto.Activate()
$wsProducto.Range("F63").Select()

# "VAB": the tab left active/selected, with B3 selected
$wsVAB.Activate()
$wsVAB.Range("B3").Select()
